$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.530.66'
$ws.Range("E2").Value = '  +0.55%  '

# Row 3
$ws.Range("D3").Value = '2.480.28'
$ws.Range("E3").Value = '  +0.73%  '

# Row 4
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.28'
$ws.Range("E5").Value = '  +0.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.91'
$ws.Range("E6").Value = '  -1.35%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.546'
$ws.Range("E7").Value = '  -0.91%  '

# Row 8
$ws.Range("E8").Value = '  -0.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.508'
$ws.Range("E9").Value = '  +2.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.79'
$ws.Range("E10").Value = '  -2.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  +0.91%  '

# Row 12
$ws.Range("E12").Value = '  +2.35%  '

# Row 13
$ws.Range("D13").Value = '2.859.96'
$ws.Range("E13").Value = '  +0.66%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.88'
$ws.Range("E14").Value = '  -1.75%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.27'
$ws.Range("E15").Value = '  +9.60%  '

# Row 16
$ws.Range("D16").Value = '2.482.01'
$ws.Range("E16").Value = '  +1.76%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.771'
$ws.Range("E17").Value = '  -1.80%  '

# Row 18
$ws.Range("D18").Value = '41.524.07'
$ws.Range("E18").Value = '  +0.63%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.59'
$ws.Range("E19").Value = '  +4.57%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0948'
$ws.Range("E20").Value = '  +2.87%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.37'
$ws.Range("E21").Value = '  +5.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.28'
$ws.Range("E22").Value = '  +0.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.57'
$ws.Range("E23").Value = '  -0.24%  '

# Row 24
$ws.Range("E24").Value = '  -1.58%  '

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("E26").Value = '  -0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.85'
$ws.Range("E27").Value = '  +3.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  +0.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.68'
$ws.Range("E29").Value = '  +0.76%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.89'
$ws.Range("E30").Value = '  -2.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.13'
$ws.Range("E31").Value = '  +3.93%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.46'
$ws.Range("E32").Value = '  -0.64%  '

# Row 33
$ws.Range("E33").Value = '  -0.65%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0758'
$ws.Range("E34").Value = '  +1.94%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.51'
$ws.Range("E35").Value = '  +3.13%  '

# Row 36
$ws.Range("E36").Value = '  -9.39%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.106'
$ws.Range("E37").Value = '  +4.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.91'
$ws.Range("E38").Value = '  -5.06%  '

# Row 39
$ws.Range("E39").Value = '  -3.18%  '

# Row 40
$ws.Range("E40").Value = '  +0.14%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.09'
$ws.Range("E41").Value = '  -4.16%  '

# Row 42
$ws.Range("E42").Value = '  -0.32%  '

# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.33'
$ws.Range("E43").Value = '  -1.81%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.959.99'
$ws.Range("E44").Value = '  -1.22%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0284'
$ws.Range("E45").Value = '  -0.14%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  -2.75%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.99'
$ws.Range("E47").Value = '  +1.93%  '

# Row 48
$ws.Range("D48").Value = '2.717.59'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.68'
$ws.Range("E49").Value = '  +0.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.17'
$ws.Range("E50").Value = '  -1.17%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.62'
$ws.Range("E51").Value = '  -2.60%  '
